$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2..18 data: Index, Week, Total Games, Count Unpredicted, Count Correct, Percent Correct
$data = @(
    @(1,  "013",    65,   5,  41, "68%"),
    @(2,  "015",     1,   0,   1, "100%"),
    @(3,  "bowls",  43,   1,  34, "81%"),
    @(4,  "004",    67,   5,  43, "69%"),
    @(5,  "001",    94,  37,  35, "61%"),
    @(6,  "011",    65,   2,  44, "70%"),
    @(7,  "007",    55,   4,  29, "57%"),
    @(8,  "014",    10,   0,   5, "50%"),
    @(9,  "006",    51,   6,  30, "67%"),
    @(10, "003",    75,  24,  36, "71%"),
    @(11, "009",    54,   3,  40, "78%"),
    @(12, "002",    85,  30,  36, "65%"),
    @(13, "005",    59,   6,  35, "66%"),
    @(14, "008",    54,   1,  35, "66%"),
    @(15, "012",    68,   6,  48, "77%"),
    @(16, "010",    65,   7,  33, "57%"),
    @(17, "totals", 911, 137, 525, "68%")
)

$lastRow = 1 + $data.Count   # rows 2..18

# Columns B (Week) and F (Percent Correct) hold values that can look
# numeric ("013", "68%", ...). Mark them as Text first so Excel stores the
# literal string instead of silently coercing to a number/percentage, then
# put the cell style back to the default "Normal" afterwards so no stray
# number-format style is left attached to the cells (matches the original
# workbook, where these cells carry no explicit style).
$bRange = $ws.Range("B2:B$lastRow")
$fRange = $ws.Range("F2:F$lastRow")
$bRange.NumberFormat = "@"
$fRange.NumberFormat = "@"

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row = $row + 1
}

$bRange.Style = "Normal"
$fRange.Style = "Normal"
